$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("SNP_entryData")

$ws2.Range("A1").Value = "entryId"
$ws2.Range("A2").Value = 1
